$d = $word.ActiveDocument

$replacements = @(
    @("82×48=3936", "81×80=6480"),
    @("26×96=2496", "64×35=2240"),
    @("26×45=1170", "97×51=4947"),
    @("83×53=4399", "15×44=660"),
    @("43×80=3440", "28×26=728"),
    @("32×17=544", "68×32=2176"),
    @("55×60=3300", "33×66=2178"),
    @("99×22=2178", "54×64=3456"),
    @("99×60=5940", "94×39=3666"),
    @("33×50=1650", "79×38=3002"),
    @("95×46=4370", "43×88=3784"),
    @("30×25=750", "75×99=7425"),
    @("66×88=5808", "90×21=1890"),
    @("66×33=2178", "39×60=2340"),
    @("88×50=4400", "36×24=864"),
    @("14×98=1372", "28×91=2548"),
    @("83×39=3237", "91×96=8736"),
    @("49×20=980", "59×55=3245"),
    @("28×63=1764", "12×37=444"),
    @("23×60=1380", "90×34=3060"),
    @("31×74=2294", "93×81=7533"),
    @("97×79=7663", "92×99=9108"),
    @("16×85=1360", "23×66=1518"),
    @("51×15=765", "59×79=4661"),
    @("62×94=5828", "72×73=5256")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
